# Apply the cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

Set-TextCell $ws 'D2' '61.171.25'
Set-TextCell $ws 'E2' '  +0.28%  '
Set-TextCell $ws 'D3' '3.369.65'
Set-TextCell $ws 'E3' '  +2.11%  '
Set-TextCell $ws 'E4' '  -0.03%  '
Set-TextCell $ws 'D5' '572.24'
Set-TextCell $ws 'E5' '  +1.86%  '
Set-TextCell $ws 'D6' '137.35'
Set-TextCell $ws 'E6' '  +8.53%  '
Set-TextCell $ws 'D8' '3.368.77'
Set-TextCell $ws 'E8' '  +2.19%  '
Set-TextCell $ws 'E9' '  -0.20%  '
Set-TextCell $ws 'E10' '  +3.98%  '
Set-TextCell $ws 'D11' '0.124'
Set-TextCell $ws 'E11' '  +4.93%  '
Set-TextCell $ws 'E12' '  +5.09%  '
Set-TextCell $ws 'D13' '3.947.46'
Set-TextCell $ws 'E13' '  +2.04%  '
Set-TextCell $ws 'E14' '  +2.13%  '
Set-TextCell $ws 'D15' '0.0000174'
Set-TextCell $ws 'E15' '  +4.18%  '
Set-TextCell $ws 'D16' '3.367.82'
Set-TextCell $ws 'E16' '  +1.74%  '
Set-TextCell $ws 'D17' '25.21'
Set-TextCell $ws 'E17' '  +2.78%  '
Set-TextCell $ws 'D18' '61.162.90'
Set-TextCell $ws 'E18' '  +0.11%  '
Set-TextCell $ws 'E19' '  +4.61%  '
Set-TextCell $ws 'D20' '13.92'
Set-TextCell $ws 'E20' '  +3.90%  '
Set-TextCell $ws 'D21' '9.36'
Set-TextCell $ws 'E21' '  +4.22%  '
Set-TextCell $ws 'D22' '378.31'
Set-TextCell $ws 'E22' '  +7.72%  '
Set-TextCell $ws 'D23' '0.569'
Set-TextCell $ws 'E23' '  +2.79%  '
Set-TextCell $ws 'D24' '3.505.17'
Set-TextCell $ws 'E24' '  +2.02%  '
Set-TextCell $ws 'E25' '  +0.10%  '
Set-TextCell $ws 'D26' '70.64'
Set-TextCell $ws 'E26' '  +2.32%  '
Set-TextCell $ws 'D27' '0.0000120'
Set-TextCell $ws 'E27' '  +12.44%  '
Set-TextCell $ws 'D28' '1.64'
Set-TextCell $ws 'E28' '  +17.11%  '
Set-TextCell $ws 'D29' '7.74'
Set-TextCell $ws 'E29' '  +8.97%  '
Set-TextCell $ws 'D30' '0.996'
Set-TextCell $ws 'E30' '  -0.44%  '
Set-TextCell $ws 'D31' '8.23'
Set-TextCell $ws 'E31' '  +5.29%  '
Set-TextCell $ws 'E32' '  +5.34%  '
Set-TextCell $ws 'D33' '2.13'
Set-TextCell $ws 'E33' '  +1.71%  '
Set-TextCell $ws 'E34' '  -0.01%  '
Set-TextCell $ws 'D35' '3.399.72'
Set-TextCell $ws 'E35' '  +2.04%  '
Set-TextCell $ws 'D36' '23.40'
Set-TextCell $ws 'E36' '  +4.17%  '
Set-TextCell $ws 'E37' '  +7.94%  '
Set-TextCell $ws 'D38' '7.03'
Set-TextCell $ws 'E38' '  +4.21%  '
Set-TextCell $ws 'E39' '  +5.76%  '
Set-TextCell $ws 'D40' '161.28'
Set-TextCell $ws 'E40' '  +0.75%  '
Set-TextCell $ws 'D41' '0.0791'
Set-TextCell $ws 'E41' '  +4.89%  '
Set-TextCell $ws 'D42' '0.999'
Set-TextCell $ws 'E42' '  -0.21%  '
Set-TextCell $ws 'B43' 'Stacks'
Set-TextCell $ws 'C43' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws 'D43' '1.71'
Set-TextCell $ws 'E43' '  +10.54%  '
Set-TextCell $ws 'B44' 'Filecoin'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D44' '4.42'
Set-TextCell $ws 'E44' '  +2.38%  '
Set-TextCell $ws 'D46' '0.760'
Set-TextCell $ws 'E46' '  +2.68%  '
Set-TextCell $ws 'D47' '1.20'
Set-TextCell $ws 'E47' '  +8.31%  '
Set-TextCell $ws 'D48' '23.07'
Set-TextCell $ws 'E48' '  +4.48%  '
Set-TextCell $ws 'D49' '6.97'
Set-TextCell $ws 'E49' '  +4.59%  '
Set-TextCell $ws 'D50' '22.89'
Set-TextCell $ws 'E50' '  +10.18%  '
Set-TextCell $ws 'D51' '2.326.30'
Set-TextCell $ws 'E51' '  +8.10%  '

Write-Host "Done: applied 89 cell updates."
